# Update NATMI LR-pair output (Ace-Bdkrb2) with newly computed TPM values.
# - Target cluster "MuSCs" rows are no longer produced by the new run, so
#   rows 4, 7 and 10 (old layout) disappear and the sheet shrinks from
#   10 data+header rows (A1:T10) to 7 (A1:T7).
# - All numeric metric columns (G:T) are refreshed with the new values,
#   and the K:L (receptor-expressing cells / detection rate) values for the
#   remaining rows also change because the cluster-size denominators shifted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data, one row per (Sending cluster, Target cluster) combination that
# survives in the refreshed output. Ligand/Receptor symbols are unchanged
# (Ace / Bdkrb2), only the metrics change.
$rows = @(
    # Sending, Target, E, F, G, H, I, J, K, L, M, N, O, P, Q, R, S, T
    @("ECs",  "ECs",  3, 1, 18.639999,           55.919997,           0.1025047374898625, 0.1025047374898625, 3, 1,                  0.367603,           1.102809,  0.5971364972068339, 0.5971364972068339, 6.852119552396999,  61.66907597157299,  0.06120931989180255,  0.06120931989180253),
    @("ECs",  "FAPs", 3, 1, 18.639999,           55.919997,           0.1025047374898625, 0.1025047374898625, 1, 0.3333333333333333, 0.2480066666666667, 0.74402,   0.4028635027931661, 0.402863502793166,  4.62284401866,      41.60559616794,     0.04129541759805999,  0.04129541759805998),
    @("FAPs", "ECs",  3, 1, 159.9051616666667,   479.715485,          0.8793475053252791, 0.8793475053252789, 3, 1,                  0.367603,           1.102809,  0.5971364972068339, 0.5971364972068339, 58.78161714415165,  529.0345542973649,  0.5250904891575049,   0.5250904891575048),
    @("FAPs", "FAPs", 3, 1, 159.9051616666667,   479.715485,          0.8793475053252791, 0.8793475053252789, 1, 0.3333333333333333, 0.2480066666666667, 0.74402,   0.4028635027931661, 0.402863502793166,  39.65754612774445,  356.9179151497,     0.3542570161677742,   0.3542570161677741),
    @("MuSCs","ECs",  3, 1, 3.300083333333333,   9.90025,             0.01814775718485843,0.01814775718485842,3, 1,                  0.367603,           1.102809,  0.5971364972068339, 0.5971364972068339, 1.213120533583333,  10.91808480225,     0.01083668815752651,  0.01083668815752651),
    @("MuSCs","FAPs", 3, 1, 3.300083333333333,   9.90025,             0.01814775718485843,0.01814775718485842,1, 0.3333333333333333, 0.2480066666666667, 0.74402,   0.4028635027931661, 0.402863502793166,  0.8184426672222223, 7.365984005,        0.007311069027331913, 0.00731106902733191)
)

# Delete the three rows whose Target cluster was "MuSCs" (old rows 4, 7, 10),
# shrinking the table to 6 data rows + header.
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(4).Delete()

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]         # A: Sending cluster
    $ws.Cells.Item($r, 4).Value = $row[1]         # D: Target cluster
    $ws.Cells.Item($r, 5).Value = $row[2]         # E
    $ws.Cells.Item($r, 6).Value = $row[3]         # F
    $ws.Cells.Item($r, 7).Value = $row[4]         # G
    $ws.Cells.Item($r, 8).Value = $row[5]         # H
    $ws.Cells.Item($r, 9).Value = $row[6]         # I
    $ws.Cells.Item($r, 10).Value = $row[7]        # J
    $ws.Cells.Item($r, 11).Value = $row[8]        # K
    $ws.Cells.Item($r, 12).Value = $row[9]        # L
    $ws.Cells.Item($r, 13).Value = $row[10]       # M
    $ws.Cells.Item($r, 14).Value = $row[11]       # N
    $ws.Cells.Item($r, 15).Value = $row[12]       # O
    $ws.Cells.Item($r, 16).Value = $row[13]       # P
    $ws.Cells.Item($r, 17).Value = $row[14]       # Q
    $ws.Cells.Item($r, 18).Value = $row[15]       # R
    $ws.Cells.Item($r, 19).Value = $row[16]       # S
    $ws.Cells.Item($r, 20).Value = $row[17]       # T
    $r++
}
